$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows before the current row 473 (Excel-style row insert,
# which shifts row 473 and everything below it down by two rows).
$ws.Rows("473:474").Insert()

# --- New row 473 ---
$ws.Range("A473").Value = 5
$ws.Range("B473").Value = "Macroferia Regional de Talca"
$ws.Range("C473").Value = "Maule"
$ws.Range("D473").Value = 44694
$ws.Range("E473").Value = 7
$ws.Range("F473").Value = "Fruta"
$ws.Range("G473").Value = 100101
$ws.Range("H473").Value = "Berries"
$ws.Range("I473").Value = 100112025
$ws.Range("J473").Value = "Frutilla"
$ws.Range("K473").Value = "Sin especificar"
$ws.Range("L473").Value = "Primera"
$ws.Range("M473").Value = 50
$ws.Range("N473").Value = 12000
$ws.Range("O473").Value = 12000
$ws.Range("P473").Value = 12000
$ws.Range("Q473").Value = "$/bandeja 7 kilos"
$ws.Range("R473").Value = "Provincia de Melipilla"
$ws.Range("S473").Value = 1714
$ws.Range("T473").Value = 7

# --- New row 474 ---
$ws.Range("A474").Value = 5
$ws.Range("B474").Value = "Macroferia Regional de Talca"
$ws.Range("C474").Value = "Maule"
$ws.Range("D474").Value = 44694
$ws.Range("E474").Value = 7
$ws.Range("F474").Value = "Fruta"
$ws.Range("G474").Value = 100101
$ws.Range("H474").Value = "Berries"
$ws.Range("I474").Value = 100112025
$ws.Range("J474").Value = "Frutilla"
$ws.Range("K474").Value = "Sin especificar"
$ws.Range("L474").Value = "Segunda"
$ws.Range("M474").Value = 30
$ws.Range("N474").Value = 7000
$ws.Range("O474").Value = 7000
$ws.Range("P474").Value = 7000
$ws.Range("Q474").Value = "$/bandeja 7 kilos"
$ws.Range("R474").Value = "Provincia de Melipilla"
$ws.Range("S474").Value = 1000
$ws.Range("T474").Value = 7
